$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: year 2021 in O4 (same look as N4) ---
$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = 2021

# --- O5: Total row (style like N5, plus 0.0 number format) ---
$ws.Range("N5").Copy($ws.Range("O5"))
$ws.Range("O5").NumberFormat = "0.0"
$ws.Range("O5").Value = 689

# --- O6..O8, O10..O14: plain data rows (style like N10, i.e. 0.0 number format) ---
$ws.Range("N10").Copy($ws.Range("O6"))
$ws.Range("O6").Value = 94.1

$ws.Range("N10").Copy($ws.Range("O7"))
$ws.Range("O7").Value = 147.1

$ws.Range("N10").Copy($ws.Range("O8"))
$ws.Range("O8").Value = 10.1

# --- O9: "-" text (style like N9, plus 0.0 number format) ---
$ws.Range("N9").Copy($ws.Range("O9"))
$ws.Range("O9").NumberFormat = "0.0"
$ws.Range("O9").Value = "-"

$ws.Range("N10").Copy($ws.Range("O10"))
$ws.Range("O10").Value = 82.1

$ws.Range("N10").Copy($ws.Range("O11"))
$ws.Range("O11").Value = 145.3

$ws.Range("N10").Copy($ws.Range("O12"))
$ws.Range("O12").Value = 98.8

$ws.Range("N10").Copy($ws.Range("O13"))
$ws.Range("O13").Value = 98.7

$ws.Range("N10").Copy($ws.Range("O14"))
$ws.Range("O14").Value = 1.8

# --- O15: "-" text (style like N9, plus 0.0 number format; re-uses style built for O9) ---
$ws.Range("N9").Copy($ws.Range("O15"))
$ws.Range("O15").NumberFormat = "0.0"
$ws.Range("O15").Value = "-"

# --- O16: bottom row (style like N16, plus 0.0 number format) ---
$ws.Range("N16").Copy($ws.Range("O16"))
$ws.Range("O16").NumberFormat = "0.0"
$ws.Range("O16").Value = 10.9

# --- Selection moves to P5 (matches the author's recorded cursor position) ---
$ws.Range("P5").Select()
